$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.511861205101013
$ws.Range("B1").Value = 4.038992881774902
$ws.Range("C1").Value = 5.909438133239746
$ws.Range("D1").Value = 1.499313235282898
$ws.Range("E1").Value = 0.8364875912666321
